# "Room with transmitter location check"
#
# This script reproduces the authors edit:
#  1. SimulationParameters: InnerOb (C7) 0 -> 1 ; Runplottype (C15)
#     "MultiRefOffCentre" -> "LOSMovedBox" ; selection moves to C7.
#  2. ObstacleMaterial: Obstacle1 Znobrat (E3) recomputed for the new
#     obstacle/transmitter position.
#  3. Obstacles: the bounding box for Obstacle1 (xmin,xmax,ymin,ymax) is
#     shifted from (0,0.45,0,0.45) to (0.45,0.75,0.45,0.75).
#  4. ObstacleCoords: obstacles 12..22 are renumbered to the zero-padded
#     012..022 form, and every one of their x/y coordinates is shifted by
#     the same +0.3 translation applied to the obstacle box (z untouched).
#  5. The workbook ends with the Obstacles sheet active (selection H2),
#     matching the final state captured in the workbook view.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. SimulationParameters
# ---------------------------------------------------------------------
$simParams = $wb.Worksheets.Item("SimulationParameters")
$simParams.Range("C7").Value = 1
$simParams.Range("C15").Value = "LOSMovedBox"
$simParams.Range("C7").Select()

# ---------------------------------------------------------------------
# 2. ObstacleMaterial
# ---------------------------------------------------------------------
$obMaterial = $wb.Worksheets.Item("ObstacleMaterial")
$obMaterial.Range("E3").Value = "(0.8855548599197771+0.005952846860713298j)"
$obMaterial.Range("E3").Value = "(0.8855885559191065+0.003942667577037837j)"

# ---------------------------------------------------------------------
# 3. Obstacles - move Obstacle1's box by (+0.45,+0.45) in x and y
# ---------------------------------------------------------------------
$obstacles = $wb.Worksheets.Item("Obstacles")
$obstacles.Range("E2").Value = 0.45
$obstacles.Range("F2").Value = 0.75
$obstacles.Range("G2").Value = 0.45
$obstacles.Range("H2").Value = 0.75

# ---------------------------------------------------------------------
# 4. ObstacleCoords - rename Obstacle12..Obstacle22 to Obstacle012..
#    Obstacle022 and translate their triangle coordinates by +0.3 in x
#    (column B) and y (column C); z (column D) is unchanged.
# ---------------------------------------------------------------------
$obCoords = $wb.Worksheets.Item("ObstacleCoords")

function Shift-Coord($value) {
    if ($value -eq 0.45) {
        return 0.75
    } else {
        return 0.45
    }
}

for ($n = 12; $n -le 22; $n++) {
    $startRow = 37 + ($n - 12) * 3
    $label = "Obstacle0{0}" -f $n
    $obCoords.Range("A$startRow").Value = $label

    for ($r = $startRow; $r -le ($startRow + 2); $r++) {
        $bCell = $obCoords.Range("B$r")
        $bCell.Value = Shift-Coord($bCell.Value2)
        $cCell = $obCoords.Range("C$r")
        $cCell.Value = Shift-Coord($cCell.Value2)
    }
}

# ---------------------------------------------------------------------
# 5. Final selection / active sheet state
# ---------------------------------------------------------------------
$obstacles.Activate()
$obstacles.Range("H2").Select()
